$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Q0" column of values being inserted at column B, pushing the
# existing Q0..Q8 (old columns B..J/K) one column to the right (C..K/L).
# Rows 2-10 already had a full 10-quarter window (B..K), so the oldest
# quarter (old column K) falls off the right-hand edge. Rows 11-20 had
# fewer than 10 quarters, so the window simply grows by one column with
# nothing dropped.
$newValues = @{
    2  = -0.58442257821662
    3  = 1.454533757567239
    4  = -1.777394389465022
    5  = -1.722070219091221
    6  = 0.3648791949059138
    7  = -0.2352699264540507
    8  = -0.05148746350304451
    9  = -0.1333319740152609
    10 = 1.614150253737389
    11 = 0.5701030647716323
    12 = 0.2202779152847414
    13 = 0.5040960054549828
    14 = 0.420735823599318
    15 = -0.1252583916527783
    16 = 0.08824118641116785
    17 = -0.1133200159455487
    18 = 0.1743923273248104
    19 = -0.4559694969238889
    20 = 0.1808172637304477
}

# Process from the bottom row up, and for each row shift the existing
# values one column to the right before writing the new value into B,
# so we never overwrite a value before it has been copied onward.
for ($row = 20; $row -ge 2; $row--) {

    # Determine how many existing data columns (starting at B=2) are
    # populated in this row, by scanning until an empty cell is found.
    # Cap at column K (11) - rows 2-10 are already a full 10-quarter
    # window, so the source column for the shift never goes past J (10):
    # the old value in K simply falls off the tracked window.
    $lastCol = 1
    for ($col = 2; $col -le 11; $col++) {
        if ($ws.Cells.Item($row, $col).Value() -ne $null) {
            $lastCol = $col
        }
    }
    $shiftFrom = [Math]::Min($lastCol, 10)

    # Shift existing values one column to the right (iterate right-to-left
    # so each write happens before its source cell is overwritten).
    for ($col = $shiftFrom; $col -ge 2; $col--) {
        $ws.Cells.Item($row, $col + 1).Value = $ws.Cells.Item($row, $col).Value()
    }

    # Write the newly-inserted leading value into column B.
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
